$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Fix the table layout so column widths are respected exactly (tblLayout type="fixed")
$t.AllowAutoFit = $false

$row = $t.Rows(1)
$row.Cells(1).Width = 21.3
$row.Cells(2).Width = 82
$row.Cells(3).Width = 76.75
$row.Cells(4).Width = 65.75
$row.Cells(5).Width = 62.1
$row.Cells(6).Width = 117.4
$row.Cells(7).Width = 99.2
